$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Range("H1").Value = "v3.3 (Mar 6, 2019)"
